$d = $word.ActiveDocument

# Remove the two "Szamo Bau Kft." paragraphs (including their paragraph marks)
$d.Content.Find.Execute(
    "Szamo Bau Kft. - Régiposta Kft., 1077 Budapest, Király u. 97 2/2. / 2015-09-02`r",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 2)

$d.Content.Find.Execute(
    "Szamo Bau Kft. - EPC Hungary Kft., 1149 Budapest, Angol u. 7. / 2015-09-02`r",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 2)

# Update the date
$d.Content.Find.Execute(
    "Budapest, 2015-10-20",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Budapest, 2015-10-23", 2)
